$d = $word.ActiveDocument

function Find-ParagraphIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text -replace "[\r\a]+$", ""
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# 1) Remove the stray _GoBack bookmark that currently sits right after "Draft #2"
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# 2) Remove the old "scratch" paragraphs that followed the Limitation bullet list:
#    - the empty bulleted (ListParagraph) paragraph
#    - "Sadasd"
#    - an empty paragraph
#    - "(OLD DRAFT)"
#    (Keeping the very last paragraph, which used to hold "The paper focuses..."
#     and will become the new Draft #3 body text below.)
$trendReportsIdx = Find-ParagraphIndexByText $d "Trend Reports, Survey Data, Employee Reports"
$oldDraftIdx = Find-ParagraphIndexByText $d "(OLD DRAFT)"

$emptyListPara = $d.Paragraphs.Item($trendReportsIdx + 1)
$oldDraftPara = $d.Paragraphs.Item($oldDraftIdx)
$deleteRange = $d.Range($emptyListPara.Range.Start, $oldDraftPara.Range.End)
$deleteRange.Delete()

# 3) Insert a new blank paragraph and a "Draft #3" paragraph right after the
#    "Trend Reports, Survey Data, Employee Reports" bullet, ahead of what used
#    to be the "The paper focuses..." paragraph.
$trendReportsPara = $d.Paragraphs.Item($trendReportsIdx)
$insertionPoint = $d.Range($trendReportsPara.Range.End, $trendReportsPara.Range.End)
$insertionPoint.InsertAfter("`rDraft #3`r")

# 4) Replace the old draft body text with the new Draft #3 body text.
$oldBodyText = "The paper focuses on developing a predictive tool for SM Hotels. The idea for now is to foretell what would happen on the efficiency of SM Hotels services if Rate of Occupancy would increase or decline. Prediction about revenues are not to be included, the goal is to come up an unbiased prediction that would increase the information for the marketing team for the to have a prepared strategy for the following days, weeks, and months"
$newBodyText = "The project will cover the automated predictive analytics system for the SM Hotels and Conventions. The system that would be created would focus on the predicting of the months and year of the people who visit the hotel in that month and how many people visited that month or year. After predicting the months and year of the people who visited in the hotel they are already alert that in, for example, at the end of the march they need to add an employee because in this month more people are visiting in the hotel."

$bodyIdx = Find-ParagraphIndexByText $d $oldBodyText
$bodyPara = $d.Paragraphs.Item($bodyIdx)
$bodyPara.Range.Find.Execute($oldBodyText, $true, $false, $false, $false, $false, $true, 1, $false, $newBodyText, 2) | Out-Null

# 5) Re-create the _GoBack bookmark right after the leading "T" of the new body
#    paragraph (i.e. splitting "The project..." into "T" + "he project...").
$bodyPara = $d.Paragraphs.Item($bodyIdx)
$afterT = $bodyPara.Range.Start + 1
$bookmarkRange = $d.Range($afterT, $afterT)
[void]$d.Bookmarks.Add("_GoBack", $bookmarkRange)
